$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 210
$ws.Range("F3").Value = 390
$ws.Range("F4").Value = 1683
$ws.Range("F5").Value = 826
$ws.Range("F6").Value = 716
$ws.Range("F7").Value = 2704
$ws.Range("F8").Value = 1357
$ws.Range("F9").Value = 2085
$ws.Range("F10").Value = 850
$ws.Range("F11").Value = 2339
$ws.Range("F12").Value = 731
$ws.Range("F13").Value = 6716
$ws.Range("F14").Value = 132
$ws.Range("F16").Value = 1271
$ws.Range("F17").Value = 1531
$ws.Range("F18").Value = 1346
$ws.Range("F21").Value = 2638
$ws.Range("D22").Value = "丰谷路35号 上海西岸艺术中心N馆"
$ws.Range("F22").Value = 1950
$ws.Range("F24").Value = 1025
$ws.Range("F25").Value = 792
$ws.Range("F26").Value = 1122
$ws.Range("F27").Value = 261
$ws.Range("F28").Value = 5398
$ws.Range("F29").Value = 293
$ws.Range("F30").Value = 1029
$ws.Range("F31").Value = 1275
$ws.Range("F32").Value = 3778
$ws.Range("F34").Value = 1709
$ws.Range("F35").Value = 1079
$ws.Range("F36").Value = 63
$ws.Range("F38").Value = 970
$ws.Range("F39").Value = 1060
$ws.Range("F40").Value = 413
$ws.Range("F41").Value = 1776
$ws.Range("F42").Value = 48
$ws.Range("F43").Value = 110
$ws.Range("F44").Value = 919
$ws.Range("F45").Value = 1053
$ws.Range("F47").Value = 519

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 79
$ws.Range("F4").Value = 7
$ws.Range("F6").Value = 445
$ws.Range("F10").Value = 400
$ws.Range("F20").Value = 608
$ws.Range("F29").Value = 71
$ws.Range("E39").Value = "2024.09.16 15:30-09.16 17:00"

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 407
$ws.Range("F10").Value = 2839
$ws.Range("F11").Value = 320
$ws.Range("F12").Value = 590
$ws.Range("F13").Value = 700
$ws.Range("F14").Value = 1211

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 210
$ws.Range("F5").Value = 390
$ws.Range("F6").Value = 407
$ws.Range("F7").Value = 2839
$ws.Range("F8").Value = 1683
$ws.Range("F9").Value = 826
$ws.Range("F10").Value = 2704
$ws.Range("F11").Value = 320
$ws.Range("F12").Value = 1357
$ws.Range("F13").Value = 850
$ws.Range("F14").Value = 2339
$ws.Range("F15").Value = 6716
$ws.Range("F16").Value = 132
$ws.Range("F17").Value = 590
$ws.Range("F19").Value = 1271
$ws.Range("F20").Value = 700
$ws.Range("F21").Value = 1531
$ws.Range("F22").Value = 1346
$ws.Range("F24").Value = 1211
$ws.Range("F25").Value = 2638
$ws.Range("D26").Value = "丰谷路35号 上海西岸艺术中心N馆"
$ws.Range("F26").Value = 1951
$ws.Range("F28").Value = 1025
$ws.Range("F29").Value = 792
$ws.Range("F30").Value = 1122
$ws.Range("F31").Value = 261
$ws.Range("F32").Value = 5398
$ws.Range("F33").Value = 293
$ws.Range("F34").Value = 1029
$ws.Range("F35").Value = 1275
$ws.Range("F36").Value = 3778
$ws.Range("F37").Value = 1709
$ws.Range("F38").Value = 1079
$ws.Range("F39").Value = 71
$ws.Range("F40").Value = 63
$ws.Range("F41").Value = 970
$ws.Range("F42").Value = 1060
$ws.Range("F43").Value = 413
$ws.Range("F44").Value = 1776
$ws.Range("F45").Value = 48
$ws.Range("F46").Value = 919
$ws.Range("F47").Value = 1053
$ws.Range("F49").Value = 519
